# Atualização de bases das ligas, do dia: 21-04-2024 às 13:33
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Swap the full data (columns B:AC) between row pairs whose
#    match entries were reordered in the source feed (column A,
#    the sequential display index, stays put on each row).
# ---------------------------------------------------------------
function Swap-Rows($sheet, $r1, $r2) {
    $rng1 = $sheet.Range("B$r1`:AC$r1")
    $rng2 = $sheet.Range("B$r2`:AC$r2")
    $v1 = $rng1.Value()
    $v2 = $rng2.Value()
    $rng1.Value = $v2
    $rng2.Value = $v1
}

Swap-Rows $ws 102 103
Swap-Rows $ws 112 114
Swap-Rows $ws 132 133
Swap-Rows $ws 179 180

# ---------------------------------------------------------------
# 2) Refresh the odds for the upcoming fixtures (rows 186-189),
#    which now correspond to the matches previously sitting in
#    rows 189-192, with updated prices. The three stale fixtures
#    that used to occupy rows 186-188 are gone.
# ---------------------------------------------------------------
function Set-MatchRow($sheet, $row, $id, $date, $home, $away, $odds) {
    $sheet.Range("B$row").Value = $id
    $sheet.Range("E$row").Value = $date
    $sheet.Range("F$row").Value = $home
    $sheet.Range("G$row").Value = $away
    $cols = @("K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $sheet.Range("$($cols[$i])$row").Value = $odds[$i]
    }
}

Set-MatchRow $ws 186 7723562 45403.5625         "Huachipato"      "Deportes Iquique"      @(2.3, 3.3, 3,   1.909, 3.75, 3.8,  -0.5,  1.925, 1.925, 2.75, 1.975, 1.875, 0, 0, 0, 0, 0)
Set-MatchRow $ws 187 7723568 45403.66666666666  "Palestino"       "Universidad de Chile"  @(3,   3.4, 2.3, 2.8,   3.4,  2.55,  0,     2.025, 1.825, 2.25, 1.8,   2.05,  0, 0, 0, 0, 0)
Set-MatchRow $ws 188 7723567 45403.77083333334  "Union Espanola"  "Deportes Copiapo"      @(1.8, 3.6, 4.5, 1.727, 3.75, 4.75, -0.75, 1.95,  1.9,   2.75, 1.975, 1.875, 0, 0, 0, 0, 0)
Set-MatchRow $ws 189 7723565 45403.875          "Cobreloa"        "Audax Italiano"        @(2,   3.5, 3.5, 2.4,   3.3,  3.1,  -0.25, 2.025, 1.825, 2.25, 1.825, 2.025, 0, 0, 0, 0, 0)

# ---------------------------------------------------------------
# 3) The old rows 190-192 (fixtures 7723568/7723567/7723565, now
#    relocated above) are removed entirely, shrinking the sheet
#    from A1:AC192 down to A1:AC189.
# ---------------------------------------------------------------
$ws.Range("A190:AC192").EntireRow.Delete() | Out-Null
